# "continue work on creds & update template"
#
# - Add two new credential rows to the ext_cred sheet (IAU-000062 / Acorn
#   Institute: a Bachelor of Computer Science and a Bachelor of
#   Informaticianistics).
# - Touch up column widths across the workbook's sheets (as happens when a
#   user resizes/auto-fits columns while reviewing).
# - Leave the workbook with the ext_cred sheet active/selected, matching
#   where the user ended up after typing the new rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# whed_levels (sheet 1) - widen the first three columns
# ---------------------------------------------------------------------
$wsLevels = $wb.Worksheets.Item("whed_levels")
$wsLevels.Columns.Item(1).ColumnWidth = 28.4133
$wsLevels.Columns.Item(2).ColumnWidth = 18.58
$wsLevels.Columns.Item(3).ColumnWidth = 28.58

# ---------------------------------------------------------------------
# whed_inst (sheet 2) - add widths for the first two columns
# ---------------------------------------------------------------------
$wsInst = $wb.Worksheets.Item("whed_inst")
$wsInst.Columns.Item(1).ColumnWidth = 14.4133
$wsInst.Columns.Item(2).ColumnWidth = 27.7467

# ---------------------------------------------------------------------
# ext_inst (sheet 3) - set widths for the first five columns and move
# the selection from B2:D2 to C2
# ---------------------------------------------------------------------
$wsExtInst = $wb.Worksheets.Item("ext_inst")
$wsExtInst.Columns.Item(1).ColumnWidth = 12.58
$wsExtInst.Columns.Item(2).ColumnWidth = 11.7467
$wsExtInst.Columns.Item(3).ColumnWidth = 14.2467
$wsExtInst.Columns.Item(4).ColumnWidth = 12.9133
$wsExtInst.Columns.Item(5).ColumnWidth = 12.9133
$wsExtInst.Activate()
$wsExtInst.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# ext_cred (sheet 4) - the actual data edit: two new credential rows,
# plus column widths now that the sheet has real data in it.
# ---------------------------------------------------------------------
$wsCred = $wb.Worksheets.Item("ext_cred")

# Institution details for both new rows.
$wsCred.Range("A2").Value = "Blue62"
$wsCred.Range("B2").Value = "Acorn Institute"
$wsCred.Range("C2").Value = "No"
$wsCred.Range("A3").Value = "Blue62"
$wsCred.Range("B3").Value = "Acorn Institute"
$wsCred.Range("C3").Value = "No"

# Course name/code for each row.
$wsCred.Range("E2").Value = "Bachelor of Computer Science"
$wsCred.Range("F2").Value = "Red71"
$wsCred.Range("E3").Value = "Bachelor of Informaticianistics"
$wsCred.Range("F3").Value = "Orange31"

# Course level (picked from the level list) for both rows.
$wsCred.Range("D2").Value = "Bachelor"
$wsCred.Range("D3").Value = "Bachelor"

# Field of study levels for each row.
$wsCred.Range("G2").Value = "02 - Information Technology"
$wsCred.Range("H2").Value = "0201 - Computer Science"
$wsCred.Range("G3").Value = "02 - Information Technology"
$wsCred.Range("H3").Value = "0208 - Informatician Mathmetician"

$wsCred.Columns.Item(1).ColumnWidth = 11.7467
$wsCred.Columns.Item(2).ColumnWidth = 15.4133
$wsCred.Columns.Item(3).ColumnWidth = 6.9133
$wsCred.Columns.Item(4).ColumnWidth = 12.9133
$wsCred.Columns.Item(5).ColumnWidth = 34.2467
$wsCred.Columns.Item(6).ColumnWidth = 26.9133
$wsCred.Columns.Item(7).ColumnWidth = 24.2467
$wsCred.Columns.Item(8).ColumnWidth = 22.9133
$wsCred.Columns.Item(9).ColumnWidth = 15.58
$wsCred.Columns.Item(10).ColumnWidth = 15.7467
$wsCred.Columns.Item(11).ColumnWidth = 10.58
$wsCred.Columns.Item(12).ColumnWidth = 12.58
$wsCred.Columns.Item(13).ColumnWidth = 10.9133
$wsCred.Columns.Item(14).ColumnWidth = 10.7467

# ext_cred ends up the active sheet, with the selection resting just past
# the newly-entered data.
$wsCred.Activate()
$wsCred.Range("G4").Select() | Out-Null
